$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (avoid numeric auto-conversion) without leaving a
# residual style on the cells: apply a Text format, write the values,
# then reset the style back to Normal (index 0) so the saved file has no
# "s" attribute on these cells, matching the original workbook layout.
$ws.Range("B2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "44.219.79"

# Row 3
$ws.Range("D3").Value = "2.264.05"
$ws.Range("E3").Value = "  +1.56%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "320.60"
$ws.Range("E5").Value = "  -1.03%  "

# Row 6
$ws.Range("D6").Value = "102.71"
$ws.Range("E6").Value = "  +3.35%  "

# Row 7
$ws.Range("E7").Value = "  -0.27%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").Value = "0.555"
$ws.Range("E9").Value = "  -1.71%  "

# Row 10
$ws.Range("D10").Value = "37.44"
$ws.Range("E10").Value = "  +1.29%  "

# Row 11
$ws.Range("D11").Value = "0.0839"
$ws.Range("E11").Value = "  +0.97%  "

# Row 12
$ws.Range("D12").Value = "7.65"
$ws.Range("E12").Value = "  -0.27%  "

# Row 13
$ws.Range("E13").Value = "  -0.95%  "

# Row 14
$ws.Range("D14").Value = "2.599.37"
$ws.Range("E14").Value = "  +1.50%  "

# Row 15
$ws.Range("D15").Value = "0.864"
$ws.Range("E15").Value = "  -0.08%  "

# Row 16
$ws.Range("D16").Value = "14.50"
$ws.Range("E16").Value = "  +0.51%  "

# Row 17
$ws.Range("D17").Value = "2.260.54"
$ws.Range("E17").Value = "  +1.72%  "

# Row 18
$ws.Range("D18").Value = "44.073.20"
$ws.Range("E18").Value = "  +2.44%  "

# Row 19
$ws.Range("D19").Value = "13.43"
$ws.Range("E19").Value = "  -5.16%  "

# Row 20
$ws.Range("E20").Value = "  +1.98%  "

# Row 21
$ws.Range("D21").Value = "6.55"
$ws.Range("E21").Value = "  -0.10%  "

# Row 22
$ws.Range("D22").Value = "65.87"
$ws.Range("E22").Value = "  +0.94%  "

# Row 23
$ws.Range("D23").Value = "3.16"
$ws.Range("E23").Value = "  -2.29%  "

# Row 24
$ws.Range("E24").Value = "  -0.69%  "

# Row 25
$ws.Range("E25").Value = "  -3.18%  "

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("D27").Value = "10.32"
$ws.Range("E27").Value = "  +3.04%  "

# Row 28
$ws.Range("D28").Value = "39.44"
$ws.Range("E28").Value = "  +7.91%  "

# Row 29
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  -2.87%  "

# Row 30
$ws.Range("D30").Value = "6.26"
$ws.Range("E30").Value = "  -1.83%  "

# Row 31
$ws.Range("D31").Value = "162.79"
$ws.Range("E31").Value = "  +5.40%  "

# Row 32
$ws.Range("D32").Value = "20.29"
$ws.Range("E32").Value = "  -0.52%  "

# Row 33
$ws.Range("D33").Value = "0.0854"
$ws.Range("E33").Value = "  -1.98%  "

# Row 34
$ws.Range("E34").Value = "  +0.27%  "

# Row 35
$ws.Range("E35").Value = "  +10.67%  "

# Row 36
$ws.Range("D36").Value = "1.98"
$ws.Range("E36").Value = "  +3.35%  "

# Row 37
$ws.Range("D37").Value = "3.07"
$ws.Range("E37").Value = "  -6.51%  "

# Row 38
$ws.Range("E38").Value = "  -1.61%  "

# Row 39
$ws.Range("D39").Value = "16.83"
$ws.Range("E39").Value = "  +20.99%  "

# Row 40
$ws.Range("E40").Value = "  -0.05%  "

# Row 41
$ws.Range("D41").Value = "4.23"
$ws.Range("E41").Value = "  -4.65%  "

# Row 42
$ws.Range("D42").Value = "0.0318"
$ws.Range("E42").Value = "  -1.71%  "

# Row 43
$ws.Range("E43").Value = "  +0.13%  "

# Row 44
$ws.Range("D44").Value = "1.788.16"
$ws.Range("E44").Value = "  +3.72%  "

# Row 45
$ws.Range("D45").Value = "0.200"
$ws.Range("E45").Value = "  -1.60%  "

# Row 46
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "82.74"
$ws.Range("E46").Value = "  -3.04%  "

# Row 47
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").Value = "75.07"
$ws.Range("E47").Value = "  -0.31%  "

# Row 48
$ws.Range("E48").Value = "  -1.12%  "

# Row 49
$ws.Range("D49").Value = "105.18"
$ws.Range("E49").Value = "  +1.88%  "

# Row 50
$ws.Range("E50").Value = "  +6.40%  "

# Row 51
$ws.Range("D51").Value = "58.66"
$ws.Range("E51").Value = "  +0.44%  "

# Reset style back to Normal so no stray style index is left on the cells
$ws.Range("B2:E51").Style = "Normal"
